# Update H column values (2nd "LS1-GA" series) for several rows on the
# "solutions" sheet. Column I holds shared formulas (H/B-1) that Excel
# will recalculate automatically once H changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solutions")

$ws.Range("H3").Value  = 7801
$ws.Range("H4").Value  = 909851
$ws.Range("H5").Value  = 52963
$ws.Range("H7").Value  = 104344
$ws.Range("H8").Value  = 1661729
$ws.Range("H11").Value = 975472
$ws.Range("H12").Value = 1379474
$ws.Range("H14").Value = 150055

# Recalculate so the dependent I-column shared formulas pick up new values.
$excel.Calculate()

# Update the active selection on the sheet to match the saved view state.
$ws.Activate()
$ws.Range("I7").Select()
